# A new observation was logged for 2026/01/11 (time=16) that belongs right
# after the existing 2026/01/11 rows (610-613) and before the 2026/12/29
# block (previously row 614). Insert a new row at 614, pushing everything
# from the old row 614 onward down by one (old 655 -> new 656), and fill
# in the new row's data.
#
# We copy the existing row 613 (same date/weekday text "2026/01/11"/"日")
# and use Insert() on the copied range so the new row 614 inherits the
# exact same cell types/formatting (inline/shared text, not an auto-
# converted date serial) instead of assigning "2026/01/11" via .Value,
# which Excel would otherwise reinterpret as a date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(613).Copy()
$ws.Rows.Item(614).Insert()
$excel.CutCopyMode = $false

# Fix up the one column that differs from the copied template row (613
# had time=13; the new row needs time=16). Date (A), weekday (B) and
# ranking (D) already match after the copy.
$ws.Range("C614").Value = 16
